# Error Calculations and Plots
# Remove the data rows for "RM 232" and "SC 92" from the sheet, shifting
# subsequent rows up (matching the source diff which drops these two rows
# and renumbers the remaining rows, shrinking the used range from
# A1:F35 to A1:F33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row whose ID (column A) is "RM 232" (currently row 26).
# Deleting this row first shifts "SC 92" up from row 28 to row 27.
$ws.Rows.Item(26).EntireRow.Delete()

# Delete the row whose ID (column A) is now "SC 92" (row 27 after the
# previous deletion).
$ws.Rows.Item(27).EntireRow.Delete()
